$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet used to hold a "CardID" header in A1 followed by a single UUID
# in A2. The card-data writer now overwrites the sheet with just the raw
# card UUIDs (no header), so clear the old header cell and (re)write the
# UUID values starting at row 2.
$ws.Range("A1").ClearContents()

$ws.Range("A2").Value = "b2175ce7-ac92-4073-bfd7-6eb30f992353"
$ws.Range("A3").Value = "e6c063f2-9c99-4f42-ad74-574fbacd06a3"
$ws.Range("A4").Value = "c2107520-2f9f-4ce4-b67d-211e61b3357c"
$ws.Range("A5").Value = "868b681d-92b0-4c69-931d-6b0a0dcbce2c"
